# Dashboard_2026.xlsx update — Janeiro 2026 numbers
# (fix: Corrigir cálculo de poupança em sync_excel.py e regenerar relatórios)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Dashboard"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Dashboard")

# Timestamp banner
$ws.Range("A2").Value = "Atualizado: 25/01/2026 17:00"

# RESUMO DO MES
$ws.Range("B7").Value = 20700        # Gastos Variaveis - Projetado
$ws.Range("C7").Value = 21602.84     # Gastos Variaveis - Real
# D7 holds a literal text percentage ("-65%" -> "4%"); pre-format the cell
# as Text so Excel stores the literal string instead of re-parsing it as a
# percentage number.
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4%"

$ws.Range("C9").Value = 22230.26     # Obra - Real

# GASTOS POR CATEGORIA
$ws.Range("C14").Value = 2942.2      # Alimentacao - Real
$ws.Range("D14").Value = 84          # Alimentacao - %

$ws.Range("C15").Value = 3166.83     # Transporte - Real
$ws.Range("D15").Value = 126         # Transporte - %

$ws.Range("C16").Value = 393         # Saude - Real
$ws.Range("D16").Value = 78          # Saude - %

$ws.Range("C17").Value = 3685.44     # Assinaturas - Real
$ws.Range("D17").Value = 92          # Assinaturas - %

$ws.Range("C18").Value = 4624.690000000001   # Compras - Real
$ws.Range("D18").Value = 115                 # Compras - %

$ws.Range("C19").Value = 4021.1      # Lazer - Real
$ws.Range("D19").Value = 114         # Lazer - %

$ws.Range("C20").Value = 1563.33     # Educacao - Real
$ws.Range("D20").Value = 104         # Educacao - %

$ws.Range("C21").Value = 163.84      # Casa - Real
$ws.Range("D21").Value = 81          # Casa - %

$ws.Range("B22").Value = 1000        # Taxas - Budget
$ws.Range("C22").Value = 1042.41     # Taxas - Real
$ws.Range("D22").Value = 104         # Taxas - %

# ------------------------------------------------------------------
# Sheet "Mensal"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mensal")

$ws.Range("B4").Value = 2942.2       # stray header total
$ws.Range("B5").Value = 3166.83      # Alimentacao
$ws.Range("B6").Value = 393          # Transporte
$ws.Range("B7").Value = 3685.44      # Saude
$ws.Range("B8").Value = 4624.690000000001  # Assinaturas
$ws.Range("B9").Value = 4021.1       # Compras
$ws.Range("B10").Value = 1563.33     # Lazer
$ws.Range("B11").Value = 163.84      # Educacao
$ws.Range("B12").Value = 1042.41     # Casa

# ------------------------------------------------------------------
# Sheet "Categorias"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Categorias")

$ws.Range("C4").Value = 2942.2
$ws.Range("D4").Value = 557.8000000000002
$ws.Range("E4").Value = 0.8406285714285714

$ws.Range("C5").Value = 3166.83
$ws.Range("D5").Value = -666.8299999999999
$ws.Range("E5").Value = 1.266732

$ws.Range("C6").Value = 393
$ws.Range("D6").Value = 107
$ws.Range("E6").Value = 0.786

$ws.Range("C7").Value = 3685.44
$ws.Range("D7").Value = 314.5599999999999
$ws.Range("E7").Value = 0.9213600000000001

$ws.Range("C8").Value = 4624.690000000001
$ws.Range("D8").Value = -624.6900000000005
$ws.Range("E8").Value = 1.1561725

$ws.Range("C9").Value = 4021.1
$ws.Range("D9").Value = -521.0999999999999
$ws.Range("E9").Value = 1.148885714285714

$ws.Range("C10").Value = 1563.33
$ws.Range("D10").Value = -63.33000000000015
$ws.Range("E10").Value = 1.04222

$ws.Range("C11").Value = 163.84
$ws.Range("D11").Value = 36.16
$ws.Range("E11").Value = 0.8192

$ws.Range("C12").Value = 1042.41
$ws.Range("D12").Value = -42.41000000000008
$ws.Range("E12").Value = 1.04241

# ------------------------------------------------------------------
# Sheet "Dados"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Dados")

$ws.Range("B3").Value = "2026-01-25T17:00:59.990258"

$ws.Range("D8").Value = 2942.2       # alimentacao
$ws.Range("D9").Value = 3166.83      # transporte
$ws.Range("D10").Value = 393         # saude
$ws.Range("D11").Value = 3685.44     # assinaturas
$ws.Range("D12").Value = 4624.690000000001   # compras
$ws.Range("D13").Value = 4021.1      # lazer
$ws.Range("D14").Value = 1563.33     # educacao
$ws.Range("D15").Value = 163.84      # casa
$ws.Range("D16").Value = 1042.41     # taxas
